$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2 through 23 contain the per-game data rows.
# Column E = birth_year, should become 1984 (was 1985).
# Column G = age, should be incremented by 1 (each player-season is a year older).
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 5).Value = 1984
    $currentAge = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 7).Value = $currentAge + 1
}
